# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-29 14:13:20
# The "Recorded By" column (G) lists the users who recorded a session, e.g. "System, dnasr281@gmail.com".
# Upstream re-ordered these so the actual user is listed before "System" (cosmetic re-sort of the
# recorded-by list), e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".
# Only the specific "System, <user>" combinations below are affected; other values
# (a lone "System", "System, admin@admin.com", or values that already start with a user) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "System, backup@backdoor.com, system" = "backup@backdoor.com, System, system"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($null -ne $current -and $replacements.ContainsKey($current)) {
        $cell.Value = $replacements[$current]
    }
}
